# Update the "修改时间" (last modified time) column on each portfolio sheet
# from 202509211502 -> 202509211517.
#
# The timestamp is stored as text (it is all digits, so it must be entered
# with a leading apostrophe to keep Excel from reinterpreting it as a
# number).

$wb = $excel.ActiveWorkbook

$newTimestamp = "'202509211517"

# Sheet "大智投资组合": rows 2-9, timestamp stored in column E
$ws1 = $wb.Worksheets.Item("大智投资组合")
$ws1.Range("E2:E9").Value = $newTimestamp

# Sheet "大成投资组合": rows 2-11, timestamp stored in column E
$ws2 = $wb.Worksheets.Item("大成投资组合")
$ws2.Range("E2:E11").Value = $newTimestamp

# Sheet "我的投资组合": rows 2-13, timestamp stored in column G
$ws3 = $wb.Worksheets.Item("我的投资组合")
$ws3.Range("G2:G13").Value = $newTimestamp
